$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (before row 422),
# which shifts the existing rows 422:449 down to 424:451 (carrying their
# values/styles along), matching the diff.
$ws.Rows("422:423").Insert()

# Populate the two newly inserted rows with this week's new data.
$ws.Cells.Item(422, 1).Value = 4
$ws.Cells.Item(422, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(422, 3).Value = "Los Lagos"
$ws.Cells.Item(422, 4).Value = 44706
$ws.Cells.Item(422, 5).Value = 10
$ws.Cells.Item(422, 6).Value = "Fruta"
$ws.Cells.Item(422, 7).Value = 100106
$ws.Cells.Item(422, 8).Value = "Oleaginosos"
$ws.Cells.Item(422, 9).Value = 100106002
$ws.Cells.Item(422, 10).Value = "Palta"
$ws.Cells.Item(422, 11).Value = "Hass"
$ws.Cells.Item(422, 12).Value = "Primera"
$ws.Cells.Item(422, 13).Value = 100
$ws.Cells.Item(422, 14).Value = 4300
$ws.Cells.Item(422, 15).Value = 4300
$ws.Cells.Item(422, 16).Value = 4300
$ws.Cells.Item(422, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(422, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(422, 19).Value = 4300
$ws.Cells.Item(422, 20).Value = 1

$ws.Cells.Item(423, 1).Value = 4
$ws.Cells.Item(423, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(423, 3).Value = "Los Lagos"
$ws.Cells.Item(423, 4).Value = 44706
$ws.Cells.Item(423, 5).Value = 10
$ws.Cells.Item(423, 6).Value = "Fruta"
$ws.Cells.Item(423, 7).Value = 100106
$ws.Cells.Item(423, 8).Value = "Oleaginosos"
$ws.Cells.Item(423, 9).Value = 100106002
$ws.Cells.Item(423, 10).Value = "Palta"
$ws.Cells.Item(423, 11).Value = "Hass"
$ws.Cells.Item(423, 12).Value = "Segunda"
$ws.Cells.Item(423, 13).Value = 50
$ws.Cells.Item(423, 14).Value = 4000
$ws.Cells.Item(423, 15).Value = 4000
$ws.Cells.Item(423, 16).Value = 4000
$ws.Cells.Item(423, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(423, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(423, 19).Value = 4000
$ws.Cells.Item(423, 20).Value = 1
